$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new blank data row at row 7 (for "DEPOVIT ...") ---
# Push current rows 7..16 (DOLO-D .. footer) down by one.
$ws.Range("A7:N7").Insert(-4121)
# Restore the row formatting (style ids) to match the data-row pattern used
# by every other item row, by copying formats from the row directly above.
$ws.Range("A6:N6").Copy()
$ws.Range("A7:N7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: insert a new blank data row at row 16 (for "سرنجات 3 سم") ---
# At this point row 16 still holds the WATER FOR INJECTION item (old row 15),
# and row 15 is the last "real" item row directly above it.
$ws.Range("A16:N16").Insert(-4121)
$ws.Range("A15:N15").Copy()
$ws.Range("A16:N16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merges for the two freshly-inserted rows (Insert() drops them)
$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("B16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()

# --- Step 3: (re)write every data row's content, rows 4-16 ---
$rows = @(
  @{ R=4;  A=1;  B="AVIVAVASC 5/160MG 28 F.C. TAB.";           H="0:0";    L=120.28; N="1:0" },
  @{ R=5;  A=2;  B="CARBAMIDE 10% CREAM 30 GM";                 H="2:0";    L=40;     N="1:0" },
  @{ R=6;  A=3;  B="CONTAFEVER N 200MG/5ML SUSP. 120ML";        H="8:0";    L=66;     N="2:0" },
  @{ R=7;  A=4;  B="DEPOVIT B12-1000MCG/ML 5 I.M. AMP";         H="2:2";    L=85;     N="1:0" },
  @{ R=8;  A=5;  B="DOLO-D PLUS ORAL SUSP. 115 ML";             H="3:0";    L=41;     N="1:0" },
  @{ R=9;  A=6;  B="EXOSIRYLIC 500 MG 20 F.C.TABS.";            H="0:1";    L=194;    N="1:0" },
  @{ R=10; A=7;  B="GLUCOVANCE 500/5MG 30 F.C.TAB.";            H="0:0";    L=74;     N="1:0" },
  @{ R=11; A=8;  B="HIBIOTIC N 600MG SUSP. 80 ML";               H="1:0";    L=92;     N="1:0" },
  @{ R=12; A=9;  B="MEGALASE SYRUP 125 ML";                      H="2:0";    L=31;     N="1:0" },
  @{ R=13; A=10; B="MINALAX 10 TABLETS";                         H="7:0";    L=36;     N="1:0" },
  @{ R=14; A=11; B="OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML";  H="0:0";    L=48;     N="2:0" },
  @{ R=15; A=12; B="WATER FOR INJECTION AMP. 5 ML";              H="7762:0"; L=5;      N="2:0" },
  @{ R=16; A=13; B="سرنجات 3 سم";                                H="-5:0";   L=10;     N="5:0" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 14).Value = $row.N
}

# --- Step 4: update the totals row (was row 15, now row 17) ---
$ws.Range("K17").Value = 842.28

Write-Output "done"
